$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the duplicated parameter row "COMP_WEIGHT_MEAN" (row 36), shifting the
# rows below it up by one.
$ws.Rows.Item(36).Delete()

# After deleting a row, Excel leaves the whole row selected where the deletion
# happened.
$ws.Activate()
$ws.Range("A36:XFD36").Select()
